# Update "US crude oil imports" workbook: revise 2015-2019 figures and add
# a new 2020 row to Table1 (Sheet1), then mirror the unpivoted rows into the
# Power-Query result table (pq sheet / Table1_2), and fix up the
# ExternalData_1 defined name so it still spans the full query result.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("pq")

# ---- 1. Revise existing 2015-2019 figures on Sheet1 (Table1) ----
# Year, Canadian imports (MMb/d), ROW imports (MMb/d)  [US exports column D unchanged]
$revisions = @(
    @(2, 3.04, 4.3230000000000004),
    @(3, 3.09, 4.76),
    @(4, 3.31, 4.6590000000000007),
    @(5, 3.62, 4.1479999999999997),
    @(6, 3.77, 3.0310000000000001)
)
foreach ($rev in $revisions) {
    $row = $rev[0]
    $ws1.Cells.Item($row, 2).Value = $rev[1]
    $ws1.Cells.Item($row, 3).Value = $rev[2]
}

# ---- 2. Append the new 2020 row to Table1 ----
$tbl1 = $ws1.ListObjects.Item("Table1")
$tbl1.ListRows.Add() | Out-Null
$ws1.Range("A7").Value = 2020
$ws1.Range("B7").Value = 3.34
$ws1.Range("C7").Value = 2.5369999999999999
$ws1.Range("D7").Value = 3.1749999999999998
$ws1.Range("C7:D7").NumberFormat = "0.0"

# ---- 3. Update the revised values inside the query-result table (pq sheet) ----
# Rows 2-16 hold the unpivoted Year/Attribute/Units/Value records produced by
# the "Table1" query for 2015-2019; keep them in sync with the Table1 edits.
$pqRevisions = @(
    @(2, 3.04),                     # 2015 Canadian imports
    @(3, 4.3230000000000004),       # 2015 ROW imports
    @(5, 3.09),                     # 2016 Canadian imports
    @(6, 4.76),                     # 2016 ROW imports
    @(8, 3.31),                     # 2017 Canadian imports
    @(9, 4.6590000000000007),       # 2017 ROW imports
    @(11, 3.62),                    # 2018 Canadian imports
    @(12, 4.1479999999999997),      # 2018 ROW imports
    @(14, 3.77),                    # 2019 Canadian imports
    @(15, 3.0310000000000001)       # 2019 ROW imports
)
foreach ($rev in $pqRevisions) {
    $ws2.Cells.Item($rev[0], 4).Value = $rev[1]
}

# ---- 4. Append the new 2020 rows (one per attribute) to the query table ----
$tbl2 = $ws2.ListObjects.Item(1)
$tbl2.ListRows.Add() | Out-Null
$tbl2.ListRows.Add() | Out-Null
$tbl2.ListRows.Add() | Out-Null

$ws2.Range("A17").Value = 2020
$ws2.Range("B17").Value = "Canadian imports "
$ws2.Range("C17").Value = "MMb/d"
$ws2.Range("D17").Value = 3.34

$ws2.Range("A18").Value = 2020
$ws2.Range("B18").Value = "ROW imports "
$ws2.Range("C18").Value = "MMb/d"
$ws2.Range("D18").Value = 2.5369999999999999

$ws2.Range("A19").Value = 2020
$ws2.Range("B19").Value = "U.S crude oil exports "
$ws2.Range("C19").Value = "MMb/d"
$ws2.Range("D19").Value = 3.1749999999999998

# ---- 5. Extend the ExternalData_1 defined name to cover the new rows ----
$wb.Names.Item("pq!ExternalData_1").RefersTo = "=pq!`$A`$1:`$D`$19"

# ---- 6. Restore the active-cell selections recorded by Excel after the edit ----
$ws1.Activate()
$ws1.Range("C11").Select()
$ws2.Activate()
$ws2.Range("E3").Select()
